$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated transition-probability matrix cells after adding more simulated games
# (values recomputed from updated underlying game counts; see commit message)
    $ws.Range("B2").Value = 0.217948717948718
    $ws.Range("C2").Value = 0.5160256410256411
    $ws.Range("J2").Value = 0.01923076923076923
    $ws.Range("P2").Value = 0.1538461538461539
    $ws.Range("S2").Value = 0.09294871794871795
    $ws.Range("B3").Value = 0.01212121212121212
    $ws.Range("C3").Value = 0.01818181818181818
    $ws.Range("J3").Value = 0.01818181818181818
    $ws.Range("P3").Value = 0.7818181818181819
    $ws.Range("S3").Value = 0.1696969696969697
    $ws.Range("J4").Value = 0.08695652173913043
    $ws.Range("P4").Value = 0.6521739130434783
    $ws.Range("S4").Value = 0.2608695652173913
    $ws.Range("B6").Value = 0.04225352112676056
    $ws.Range("D6").Value = 0.0187793427230047
    $ws.Range("F6").Value = 0.04694835680751173
    $ws.Range("J6").Value = 0.2676056338028169
    $ws.Range("O6").Value = 0.03286384976525822
    $ws.Range("Q6").Value = 0.2488262910798122
    $ws.Range("R6").Value = 0.05633802816901409
    $ws.Range("S6").Value = 0.2863849765258216
    $ws.Range("B7").Value = 0.1075268817204301
    $ws.Range("D7").Value = 0.02150537634408602
    $ws.Range("F7").Value = 0.04838709677419355
    $ws.Range("J7").Value = 0.1182795698924731
    $ws.Range("O7").Value = 0.02150537634408602
    $ws.Range("Q7").Value = 0.1559139784946237
    $ws.Range("R7").Value = 0.1021505376344086
    $ws.Range("S7").Value = 0.4247311827956989
    $ws.Range("B8").Value = 0.09490740740740741
    $ws.Range("D8").Value = 0.01388888888888889
    $ws.Range("F8").Value = 0.05787037037037037
    $ws.Range("J8").Value = 0.1157407407407407
    $ws.Range("O8").Value = 0.01157407407407407
    $ws.Range("Q8").Value = 0.2337962962962963
    $ws.Range("R8").Value = 0.05555555555555555
    $ws.Range("S8").Value = 0.4166666666666667
    $ws.Range("B9").Value = 0.08900523560209424
    $ws.Range("D9").Value = 0.02617801047120419
    $ws.Range("F9").Value = 0.05235602094240838
    $ws.Range("J9").Value = 0.07853403141361257
    $ws.Range("O9").Value = 0.02617801047120419
    $ws.Range("Q9").Value = 0.2408376963350785
    $ws.Range("R9").Value = 0.1047120418848168
    $ws.Range("S9").Value = 0.3821989528795812
    $ws.Range("B10").Value = 0.1123188405797101
    $ws.Range("D10").Value = 0.02246376811594203
    $ws.Range("E10").Value = 0.0007246376811594203
    $ws.Range("F10").Value = 0.06086956521739131
    $ws.Range("J10").Value = 0.1130434782608696
    $ws.Range("O10").Value = 0.02391304347826087
    $ws.Range("Q10").Value = 0.2253623188405797
    $ws.Range("R10").Value = 0.07101449275362319
    $ws.Range("S10").Value = 0.3702898550724638
    $ws.Range("G11").Value = 0.1694915254237288
    $ws.Range("J11").Value = 0.08813559322033898
    $ws.Range("K11").Value = 0.223728813559322
    $ws.Range("L11").Value = 0.5084745762711864
    $ws.Range("S11").Value = 0.01016949152542373
    $ws.Range("G12").Value = 0.7044025157232704
    $ws.Range("J12").Value = 0.2138364779874214
    $ws.Range("K12").Value = 0.01257861635220126
    $ws.Range("L12").Value = 0.03144654088050314
    $ws.Range("S12").Value = 0.03773584905660377
    $ws.Range("G13").Value = 0.6097560975609756
    $ws.Range("J13").Value = 0.3658536585365854
    $ws.Range("S13").Value = 0.02439024390243903
    $ws.Range("F15").Value = 0.02242152466367713
    $ws.Range("H15").Value = 0.1300448430493273
    $ws.Range("I15").Value = 0.05829596412556054
    $ws.Range("J15").Value = 0.3452914798206278
    $ws.Range("K15").Value = 0.04932735426008968
    $ws.Range("M15").Value = 0.0179372197309417
    $ws.Range("O15").Value = 0.04932735426008968
    $ws.Range("S15").Value = 0.3273542600896861
    $ws.Range("H16").Value = 0.185
    $ws.Range("I16").Value = 0.09
    $ws.Range("J16").Value = 0.44
    $ws.Range("K16").Value = 0.09
    $ws.Range("M16").Value = 0.005
    $ws.Range("N16").Value = 0.005
    $ws.Range("O16").Value = 0.055
    $ws.Range("S16").Value = 0.105
    $ws.Range("F17").Value = 0.01879699248120301
    $ws.Range("H17").Value = 0.1541353383458647
    $ws.Range("I17").Value = 0.07142857142857142
    $ws.Range("J17").Value = 0.4530075187969925
    $ws.Range("K17").Value = 0.09962406015037593
    $ws.Range("M17").Value = 0.01691729323308271
    $ws.Range("O17").Value = 0.07706766917293233
    $ws.Range("S17").Value = 0.1090225563909774
    $ws.Range("F18").Value = 0.02906976744186046
    $ws.Range("H18").Value = 0.1802325581395349
    $ws.Range("I18").Value = 0.0755813953488372
    $ws.Range("J18").Value = 0.436046511627907
    $ws.Range("K18").Value = 0.1046511627906977
    $ws.Range("M18").Value = 0.02325581395348837
    $ws.Range("O18").Value = 0.05232558139534884
    $ws.Range("S18").Value = 0.09883720930232558
    $ws.Range("F19").Value = 0.01630434782608696
    $ws.Range("H19").Value = 0.1940993788819876
    $ws.Range("I19").Value = 0.08462732919254658
    $ws.Range("J19").Value = 0.4052795031055901
    $ws.Range("K19").Value = 0.09627329192546584
    $ws.Range("M19").Value = 0.01940993788819876
    $ws.Range("N19").Value = 0.0007763975155279503
    $ws.Range("O19").Value = 0.06133540372670807
    $ws.Range("S19").Value = 0.1218944099378882
